# Auto-generated: update cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.942.64'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.647.84'
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").Value = '  +0.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5115'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.95%  '

$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2587'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06435'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07783'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.325'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.97%  '

$ws.Range("D13").Value = '1.658.23'
$ws.Range("E13").Value = '  +1.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5481'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").Value = '0.0₅7909'
$ws.Range("E15").Value = '  -0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.46%  '

$ws.Range("D17").Value = '26.036.11'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '199.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.465'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.080'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.008'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.866'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1152'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.913'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.75%  '

$ws.Range("E29").Value = '  +0.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05044'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.297'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.211'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.548'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.366'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8968'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.593'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.88%  '

$ws.Range("D37").Value = '1.140.37'
$ws.Range("E37").Value = '  -2.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5564'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01566'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.007'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.685'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8187'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("E44").Value = '  +9.37%  '

$ws.Range("D45").Value = '1.787.06'
$ws.Range("E45").Value = '  +0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4537'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05095'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.09592'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '

